$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D12").Value = "13:30-zeit"
$ws.Range("A12").Value = "27.12.2023"
$ws.Range("D12").NumberFormat = "h:mm"

$ws.Range("A12").Select()
